$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("E3").Value = "2016-03-25 07:34:35"
$wsZh.Range("H3").Value = "2016-03-25 07:35:09"

$wsDe.Range("E3").Value = "2016-03-25 07:34:39"
$wsDe.Range("H3").Value = "2016-03-25 07:35:17"
